$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'304.01"
$ws.Range("E2").Value = "'4.02%"
$ws.Range("D3").Value = "'35.65"
$ws.Range("E3").Value = "'14.81%"
$ws.Range("D4").Value = "'5.098"
$ws.Range("E4").Value = "'2.48%"
$ws.Range("D5").Value = "'0.07803"
$ws.Range("E5").Value = "'4.64%"
$ws.Range("D6").Value = "'2.265"
$ws.Range("E6").Value = "'0.98%"
$ws.Range("D7").Value = "'8.095"
$ws.Range("E7").Value = "'4.26%"
$ws.Range("D8").Value = "'4.002"
$ws.Range("E8").Value = "'6.06%"
$ws.Range("D9").Value = "'0.9296"
$ws.Range("D10").Value = "'0.09711"
$ws.Range("E10").Value = "'3.30%"
$ws.Range("D11").Value = "'0.1821"
$ws.Range("E11").Value = "'4.99%"
$ws.Range("D12").Value = "'0.08743"
$ws.Range("E12").Value = "'5.37%"
$ws.Range("E13").Value = "'4.52%"
$ws.Range("D14").Value = "'0.09929"
$ws.Range("E14").Value = "'0.12%"
$ws.Range("D15").Value = "'0.001480"
$ws.Range("E15").Value = "'-1.52%"
$ws.Range("D16").Value = "'0.005775"
$ws.Range("E16").Value = "'0.92%"
$ws.Range("D17").Value = "'3.489"
$ws.Range("E17").Value = "'0.62%"
$ws.Range("D18").Value = "'2.127"
$ws.Range("E18").Value = "'-3.17%"
$ws.Range("E19").Value = "'3.05%"
$ws.Range("E20").Value = "'1.02%"
$ws.Range("D21").Value = "'4.561"
$ws.Range("E21").Value = "'11.33%"
$ws.Range("E22").Value = "'5.33%"
$ws.Range("D23").Value = "'0.04675"
$ws.Range("E23").Value = "'3.12%"
$ws.Range("E24").Value = "'1.65%"
$ws.Range("D25").Value = "'0.004491"
$ws.Range("E25").Value = "'5.38%"
$ws.Range("D26").Value = "'0.0001302"
$ws.Range("E26").Value = "'0.18%"
$ws.Range("D27").Value = "'0.0002703"
$ws.Range("E27").Value = "'-20.36%"
$ws.Range("D39").Value = "'0.01754"
$ws.Range("E39").Value = "'8.26%"
$ws.Range("D40").Value = "'0.04704"
$ws.Range("D41").Value = "'0.007831"
$ws.Range("E41").Value = "'5.07%"
$ws.Range("D42").Value = "'0.1418"
$ws.Range("E42").Value = "'4.51%"
$ws.Range("D43").Value = "'0.008581"
$ws.Range("E43").Value = "'-12.85%"
$ws.Range("D44").Value = "'0.002294"
$ws.Range("E44").Value = "'3.34%"
$ws.Range("D45").Value = "'0.009220"
$ws.Range("E45").Value = "'2.99%"
$ws.Range("D46").Value = "'0.00006138"
$ws.Range("E46").Value = "'0.57%"
$ws.Range("E47").Value = "'0.17%"
$ws.Range("D48").Value = "'3.943"
$ws.Range("E48").Value = "'41.08%"
$ws.Range("D49").Value = "'0.002694"
$ws.Range("E49").Value = "'34.72%"
$ws.Range("D50").Value = "'0.00002103"
$ws.Range("E50").Value = "'0.17%"
$ws.Range("D51").Value = "'0.0002003"
$ws.Range("E51").Value = "'0.17%"
